# Adds "SamplesTab" and "FilesTab" rows (rows 3 and 4) to the startup sheet,
# mirroring the structure of the existing "CasesTab" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tab-name labels first, so they claim the lower shared-string indices
# (matches the order the rows were authored in).
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"

# --- Row 3: SamplesTab ---
$ws.Range("B3").Value = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE s.study_acronym IN ["C"]  
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`
'@
# StatQuery (col C), Neo4j output file (col D) and Web output file (col E)
# are identical across tabs, so copy them straight from the CasesTab row.
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2

# --- Row 4: FilesTab ---
$ws.Range("B4").Value = @'
MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE s.study_acronym IN ["C"]  
WITH
        f, parent,p, ss, d,tp, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`
    order by f.file_name
'@
$ws.Range("C4").Value = $ws.Range("C2").Value2
$ws.Range("D4").Value = $ws.Range("D2").Value2
$ws.Range("E4").Value = $ws.Range("E2").Value2

# Wrap the long query / stat-query text, matching columns B & C on row 2.
$ws.Range("B3:C4").WrapText = $true

# Let the new rows size themselves like the existing CasesTab row.
$ws.Rows.Item(3).RowHeight = 360
$ws.Rows.Item(4).RowHeight = 409.6

# Column widths grow slightly now that rows 3 & 4 provide new "best fit" text.
$ws.Columns.Item(1).ColumnWidth = 12.72
$ws.Columns.Item(2).ColumnWidth = 74.94
$ws.Columns.Item(3).ColumnWidth = 50.72
$ws.Columns.Item(4).ColumnWidth = 44.05
$ws.Columns.Item(5).ColumnWidth = 42.05

# Update the view: zoomed in a bit more, scrolled to the new rows, selection on them.
$ws.Application.ActiveWindow.Zoom = 60
$ws.Application.ActiveWindow.ScrollRow = 3
[void]$ws.Range("C2:E4").Select()
$ws.Application.ActiveCell = $ws.Range("C2")

Write-Output "Added SamplesTab and FilesTab rows"
